# COREESG_holdings.xlsx - "Add files via upload"
# Refreshes the model-holdings snapshot on Sheet1:
#   - rolls the "as of" date in the confidential disclaimer (A10) forward
#     from 2021-06-14 to 2021-07-07
#   - updates the Weight (D2:D6) and Percent Change (E2:E7) figures to the
#     new snapshot's numbers (percent-change column is reset to 0 for the
#     new period)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; lift protection so the cells can be written,
# then restore it once the edits are in place.
$ws.Unprotect()

# --- Confidential disclaimer date -----------------------------------
# Use Value2 (not Value) so the shared-string cell is updated in place
# rather than marshalled through the COM Variant stub.
$ws.Range("A10").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-07 for illustrative purposes only and are subject to change."

# --- NULG (row 2) -----------------------------------------------------
$ws.Range("D2").Value = 0.2640039196952417
$ws.Range("E2").Value = 0

# --- NULV (row 3) -------------------------------------------------------
$ws.Range("D3").Value = 0.5292023440231703
$ws.Range("E3").Value = 0

# --- NUMG (row 4) -------------------------------------------------------
$ws.Range("D4").Value = 0.05282862965531705
$ws.Range("E4").Value = 0

# --- NUMV (row 5) -------------------------------------------------------
$ws.Range("D5").Value = 0.09606621599336043
$ws.Range("E5").Value = 0

# --- NUSC (row 6) -------------------------------------------------------
$ws.Range("D6").Value = 0.05789889063291045
$ws.Range("E6").Value = 0

# --- Total (row 7) -------------------------------------------------------
$ws.Range("E7").Value = 0

# Restore protection (matches the shipped workbook's protected sheet).
$ws.Protect()
